$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3747116666666667
$ws.Range("H2").Value = 1.124135
$ws.Range("I2").Value = 0.3914669751594584
$ws.Range("J2").Value = 0.3914669751594584
$ws.Range("M2").Value = 1.334383666666667
$ws.Range("N2").Value = 4.003151
$ws.Range("O2").Value = 0.1312069045987744
$ws.Range("P2").Value = 0.1312069045987744
$ws.Range("Q2").Value = 0.5000091277094445
$ws.Range("R2").Value = 4.500082149385
$ws.Range("S2").Value = 0.05136317006331784
$ws.Range("T2").Value = 0.05136317006331784

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3747116666666667
$ws.Range("H3").Value = 1.124135
$ws.Range("I3").Value = 0.3914669751594584
$ws.Range("J3").Value = 0.3914669751594584
$ws.Range("N3").Value = 7.432386999999999
$ws.Range("O3").Value = 0.2436032245723858
$ws.Range("P3").Value = 0.2436032245723858
$ws.Range("Q3").Value = 0.9283340400272222
$ws.Range("R3").Value = 8.355006360245
$ws.Range("S3").Value = 0.09536261746244211
$ws.Range("T3").Value = 0.09536261746244214

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3747116666666667
$ws.Range("H4").Value = 1.124135
$ws.Range("I4").Value = 0.3914669751594584
$ws.Range("J4").Value = 0.3914669751594584
$ws.Range("M4").Value = 6.358226000000001
$ws.Range("N4").Value = 19.074678
$ws.Range("O4").Value = 0.6251898708288398
$ws.Range("P4").Value = 0.6251898708288398
$ws.Range("Q4").Value = 2.382501461503334
$ws.Range("R4").Value = 21.44251315353
$ws.Range("S4").Value = 0.2447411876336984
$ws.Range("T4").Value = 0.2447411876336984

$ws.Range("I5").Value = 0.4195022558883632
$ws.Range("J5").Value = 0.4195022558883631
$ws.Range("M5").Value = 1.334383666666667
$ws.Range("N5").Value = 4.003151
$ws.Range("O5").Value = 0.1312069045987744
$ws.Range("P5").Value = 0.1312069045987744
$ws.Range("Q5").Value = 0.535817758199
$ws.Range("R5").Value = 4.822359823791
$ws.Range("S5").Value = 0.0550415924673151
$ws.Range("T5").Value = 0.05504159246731509

$ws.Range("I6").Value = 0.4195022558883632
$ws.Range("J6").Value = 0.4195022558883631
$ws.Range("N6").Value = 7.432386999999999
$ws.Range("O6").Value = 0.2436032245723858
$ws.Range("P6").Value = 0.2436032245723858
$ws.Range("Q6").Value = 0.994817567563
$ws.Range("S6").Value = 0.1021921022497954
$ws.Range("T6").Value = 0.1021921022497954

$ws.Range("I7").Value = 0.4195022558883632
$ws.Range("J7").Value = 0.4195022558883631
$ws.Range("M7").Value = 6.358226000000001
$ws.Range("N7").Value = 19.074678
$ws.Range("O7").Value = 0.6251898708288398
$ws.Range("P7").Value = 0.6251898708288398
$ws.Range("Q7").Value = 2.553126575622001
$ws.Range("R7").Value = 22.978139180598
$ws.Range("S7").Value = 0.2622685611712527
$ws.Range("T7").Value = 0.2622685611712526

$ws.Range("G8").Value = 0.18094
$ws.Range("H8").Value = 0.54282
$ws.Range("I8").Value = 0.1890307689521785
$ws.Range("J8").Value = 0.1890307689521785
$ws.Range("M8").Value = 1.334383666666667
$ws.Range("N8").Value = 4.003151
$ws.Range("O8").Value = 0.1312069045987744
$ws.Range("P8").Value = 0.1312069045987744
$ws.Range("Q8").Value = 0.2414433806466666
$ws.Range("R8").Value = 2.17299042582
$ws.Range("S8").Value = 0.02480214206814144
$ws.Range("T8").Value = 0.02480214206814144

$ws.Range("G9").Value = 0.18094
$ws.Range("H9").Value = 0.54282
$ws.Range("I9").Value = 0.1890307689521785
$ws.Range("J9").Value = 0.1890307689521785
$ws.Range("N9").Value = 7.432386999999999
$ws.Range("O9").Value = 0.2436032245723858
$ws.Range("P9").Value = 0.2436032245723858
$ws.Range("Q9").Value = 0.4482720345933333
$ws.Range("R9").Value = 4.034448311339999
$ws.Range("S9").Value = 0.0460485048601483
$ws.Range("T9").Value = 0.04604850486014832

$ws.Range("G10").Value = 0.18094
$ws.Range("H10").Value = 0.54282
$ws.Range("I10").Value = 0.1890307689521785
$ws.Range("J10").Value = 0.1890307689521785
$ws.Range("M10").Value = 6.358226000000001
$ws.Range("N10").Value = 19.074678
$ws.Range("O10").Value = 0.6251898708288398
$ws.Range("P10").Value = 0.6251898708288398
$ws.Range("Q10").Value = 1.15045741244
$ws.Range("R10").Value = 10.35411671196
$ws.Range("S10").Value = 0.1181801220238887
$ws.Range("T10").Value = 0.1181801220238887
